# "July 28 @11:02am - updating csv file"
# The state names in column A (rows 2-49) were re-typed/re-pasted in lowercase
# (matching a refreshed source CSV), the view was scrolled/re-selected, and the
# columns were (auto-)sized to fit the new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New (lowercased) state names for A2:A49, in the same row order as before.
$states = @(
    "alabama",
    "arizona",
    "arkansas",
    "california",
    "colorado",
    "connecticut",
    "delaware",
    "florida",
    "georgia",
    "idaho",
    "illinois",
    "indiana",
    "iowa",
    "kansas",
    "kentucky",
    "louisiana",
    "maine",
    "maryland",
    "massachusetts",
    "michigan",
    "minnesota",
    "mississippi",
    "missouri",
    "montana",
    "nebraska",
    "nevada",
    "new hampshire",
    "new jersey",
    "new mexico",
    "new york",
    "north carolina",
    "north dakota",
    "ohio",
    "oklahoma",
    "oregon",
    "pennsylvania",
    "rhode island",
    "south carolina",
    "south dakota",
    "tennessee",
    "texas",
    "utah",
    "vermont",
    "virginia",
    "washington",
    "west Virginia",
    "wisconsin",
    "wyoming"
)

for ($i = 0; $i -lt $states.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $states[$i]
}

# Columns were widened to fit the refreshed data.
$ws.Columns.Item(1).ColumnWidth = 34.5
$ws.Columns.Item(2).ColumnWidth = 20
$ws.Columns.Item(3).ColumnWidth = 21.333333333333332
$ws.Columns.Item(4).ColumnWidth = 21.333333333333332

# Scroll position / selection moved.
$ws.Range("A18").Select()
$ws.Range("A50").Select()
